$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "36.00", "1.00") keep their original text representation
# instead of being auto-converted to numbers by Excel.
$touchedAddrs = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'E8', 'E9', 'E10', 'D11', 'E11', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'E19', 'E20', 'D21', 'E21', 'D22', 'E22', 'E23', 'E24', 'D25', 'E25', 'E26', 'D27', 'E27', 'D28', 'E28', 'E30', 'D31', 'E31', 'D32', 'E32', 'E33', 'D34', 'E34', 'D35', 'D36', 'E36', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'E44', 'E45', 'D46', 'E46', 'E47', 'D48', 'E48', 'B49', 'C49', 'D49', 'E49', 'B50', 'C50', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $touchedAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.171.19'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '3.736.30'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '592.37'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '167.14'
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").Value = '3.732.95'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '0.0000258'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").Value = '36.00'
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").Value = '4.361.85'
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("D16").Value = '3.729.68'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").Value = '68.097.74'
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").Value = '17.87'
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").Value = '10.65'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '465.56'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  +8.77%  '
$ws.Range("D25").Value = '83.83'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("D27").Value = '11.85'
$ws.Range("E27").Value = '  -1.12%  '
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").Value = '7.27'
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").Value = '29.74'
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("E33").Value = '  -2.94%  '
$ws.Range("D34").Value = '9.14'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("D36").Value = '3.688.91'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = '3.45'
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = '0.137'
$ws.Range("E39").Value = '  +0.78%  '
$ws.Range("D40").Value = '0.988'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").Value = '5.76'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("E44").Value = '  +15.02%  '
$ws.Range("E45").Value = '  -1.69%  '
$ws.Range("D46").Value = '46.66'
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = '8.41'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '388.09'
$ws.Range("E49").Value = '  -1.42%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '144.11'
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("D51").Value = '2.740.68'
$ws.Range("E51").Value = '  +2.35%  '
